$wb = $excel.ActiveWorkbook

# --- Sheet 1: Significant Components ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range('C2').Value = '[''QEXTRCT'' ''QEDLESHI'' ''QNOHLTH'' ''QPOVTY'' ''QESL'' ''QHISPC'' ''PPUNIT'' ''PERCAP'']'
$ws1.Range('C3').Value = '[''QPOVTY'' ''QFAM'' ''QFHH'' ''QBLACK'' ''QSERV'' ''PERCAP'' ''QRICH'']'
$ws1.Range('C4').Value = '[''QFEMLBR'' ''QFEMALE'' ''QAGEDEP'']'
$ws1.Range('C6').Value = '[''PPUNIT'' ''QRENTER'' ''QNOAUTO'']'

# --- Sheet 2: Loading Factors ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 1).Value = 'QEXTRCT'
$ws2.Cells.Item(2, 2).Value = 0.7518642484417047
$ws2.Cells.Item(2, 3).Value = 0.1080633726986945
$ws2.Cells.Item(2, 4).Value = -0.184153634128851
$ws2.Cells.Item(2, 5).Value = -0.04880562950779412
$ws2.Cells.Item(2, 6).Value = 0.02811185968668175
$ws2.Cells.Item(3, 1).Value = 'QEDLESHI'
$ws2.Cells.Item(3, 2).Value = 0.8692554802113335
$ws2.Cells.Item(3, 3).Value = 0.2820098251129958
$ws2.Cells.Item(3, 4).Value = -0.09611805635391178
$ws2.Cells.Item(3, 5).Value = -0.04435359362502852
$ws2.Cells.Item(3, 6).Value = 0.01766608572966158
$ws2.Cells.Item(4, 1).Value = 'QNOHLTH'
$ws2.Cells.Item(4, 2).Value = 0.7427426205749446
$ws2.Cells.Item(4, 3).Value = 0.3825258292209991
$ws2.Cells.Item(4, 4).Value = -0.09377154783409625
$ws2.Cells.Item(4, 5).Value = -0.1242320404236909
$ws2.Cells.Item(4, 6).Value = 0.1047037153226807
$ws2.Cells.Item(5, 1).Value = 'QPOVTY'
$ws2.Cells.Item(5, 2).Value = 0.4926703472929525
$ws2.Cells.Item(5, 3).Value = 0.4936971001788416
$ws2.Cells.Item(5, 4).Value = 0.01006171257136397
$ws2.Cells.Item(5, 5).Value = -0.139980912263577
$ws2.Cells.Item(5, 6).Value = 0.359303199780186
$ws2.Cells.Item(6, 1).Value = 'QESL'
$ws2.Cells.Item(6, 2).Value = 0.8669870460373705
$ws2.Cells.Item(6, 3).Value = 0.105588759382269
$ws2.Cells.Item(6, 4).Value = -0.1218582062889729
$ws2.Cells.Item(6, 5).Value = -0.1645406179565412
$ws2.Cells.Item(6, 6).Value = 0.1246202480462146
$ws2.Cells.Item(7, 1).Value = 'QHISPC'
$ws2.Cells.Item(7, 2).Value = 0.8803370942915406
$ws2.Cells.Item(7, 3).Value = 0.1202941229291331
$ws2.Cells.Item(7, 4).Value = -0.1034740124261152
$ws2.Cells.Item(7, 5).Value = -0.205587078981756
$ws2.Cells.Item(7, 6).Value = -0.1433395948207157
$ws2.Cells.Item(8, 1).Value = 'PPUNIT'
$ws2.Cells.Item(8, 2).Value = 0.5203431752795143
$ws2.Cells.Item(8, 3).Value = 0.1723556415683452
$ws2.Cells.Item(8, 4).Value = -0.008531183769641909
$ws2.Cells.Item(8, 5).Value = -0.0411856994804425
$ws2.Cells.Item(8, 6).Value = -0.6351977636422353
$ws2.Cells.Item(9, 1).Value = 'QFAM'
$ws2.Cells.Item(9, 2).Value = 0.1243576377146733
$ws2.Cells.Item(9, 3).Value = 0.66294007956358
$ws2.Cells.Item(9, 4).Value = 0.07601970061597395
$ws2.Cells.Item(9, 5).Value = -0.1184125638687395
$ws2.Cells.Item(9, 6).Value = 0.2379118979678956
$ws2.Cells.Item(10, 1).Value = 'QFHH'
$ws2.Cells.Item(10, 2).Value = 0.2285252645071095
$ws2.Cells.Item(10, 3).Value = 0.7131346566768045
$ws2.Cells.Item(10, 4).Value = 0.2214646517180176
$ws2.Cells.Item(10, 5).Value = -0.07783278306666057
$ws2.Cells.Item(10, 6).Value = -0.04601151712576505
$ws2.Cells.Item(11, 1).Value = 'QBLACK'
$ws2.Cells.Item(11, 2).Value = -0.2749361517443922
$ws2.Cells.Item(11, 3).Value = 0.7046485571608363
$ws2.Cells.Item(11, 4).Value = 0.05019711414732955
$ws2.Cells.Item(11, 5).Value = 0.1307361802183043
$ws2.Cells.Item(11, 6).Value = 0.1670202047959224
$ws2.Cells.Item(12, 1).Value = 'QSERV'
$ws2.Cells.Item(12, 2).Value = 0.3709861023078258
$ws2.Cells.Item(12, 3).Value = 0.5379072869833532
$ws2.Cells.Item(12, 4).Value = 0.005759199128988558
$ws2.Cells.Item(12, 5).Value = -0.08757769803385232
$ws2.Cells.Item(12, 6).Value = 0.1559839714480021
$ws2.Cells.Item(13, 1).Value = 'PERCAP'
$ws2.Cells.Item(13, 2).Value = 0.5037823282541694
$ws2.Cells.Item(13, 3).Value = 0.7061433337952319
$ws2.Cells.Item(13, 4).Value = -0.04688443481784117
$ws2.Cells.Item(13, 5).Value = -0.1072480938194731
$ws2.Cells.Item(13, 6).Value = -0.1020625141790412
$ws2.Cells.Item(14, 1).Value = 'QRICH'
$ws2.Cells.Item(14, 2).Value = 0.4038778847879076
$ws2.Cells.Item(14, 3).Value = 0.6463639839109264
$ws2.Cells.Item(14, 4).Value = -0.06158260922857307
$ws2.Cells.Item(14, 5).Value = -0.1194534731552246
$ws2.Cells.Item(14, 6).Value = 0.002273121035538215
$ws2.Cells.Item(15, 1).Value = 'QRENTER'
$ws2.Cells.Item(15, 2).Value = 0.1727227398097679
$ws2.Cells.Item(15, 3).Value = 0.3584593755415023
$ws2.Cells.Item(15, 4).Value = -0.05642808994004768
$ws2.Cells.Item(15, 5).Value = -0.4605499177939512
$ws2.Cells.Item(15, 6).Value = 0.6579700024315055
$ws2.Cells.Item(16, 1).Value = 'QNOAUTO'
$ws2.Cells.Item(16, 2).Value = 0.1891304089967761
$ws2.Cells.Item(16, 3).Value = 0.3968119693850521
$ws2.Cells.Item(16, 4).Value = -0.001312813646046432
$ws2.Cells.Item(16, 5).Value = 0.0829157104879002
$ws2.Cells.Item(16, 6).Value = 0.5521493306898304
$ws2.Cells.Item(17, 1).Value = 'QFEMLBR'
$ws2.Cells.Item(17, 2).Value = -0.4233895386433378
$ws2.Cells.Item(17, 3).Value = 0.2565876909270864
$ws2.Cells.Item(17, 4).Value = 0.5549136206110808
$ws2.Cells.Item(17, 5).Value = -0.009510215744995571
$ws2.Cells.Item(17, 6).Value = -0.03382138915842849
$ws2.Cells.Item(18, 1).Value = 'QFEMALE'
$ws2.Cells.Item(18, 2).Value = -0.1201784102251344
$ws2.Cells.Item(18, 3).Value = 0.08442001701203185
$ws2.Cells.Item(18, 4).Value = 0.9476420859302027
$ws2.Cells.Item(18, 5).Value = 0.01306748493548795
$ws2.Cells.Item(18, 6).Value = -0.03709586765110615
$ws2.Cells.Item(19, 1).Value = 'QAGEDEP'
$ws2.Cells.Item(19, 2).Value = -0.1133192053356318
$ws2.Cells.Item(19, 3).Value = -0.06739481110732691
$ws2.Cells.Item(19, 4).Value = 0.7290404171908992
$ws2.Cells.Item(19, 5).Value = 0.4777003064674182
$ws2.Cells.Item(19, 6).Value = 0.06612983910678825
$ws2.Cells.Item(20, 1).Value = 'MEDAGE'
$ws2.Cells.Item(20, 2).Value = -0.3065470993356068
$ws2.Cells.Item(20, 3).Value = -0.3702019246123122
$ws2.Cells.Item(20, 4).Value = 0.04864970172662886
$ws2.Cells.Item(20, 5).Value = 0.6451312986528182
$ws2.Cells.Item(20, 6).Value = 0.09075704017158054
$ws2.Cells.Item(21, 1).Value = 'QSSBEN'
$ws2.Cells.Item(21, 2).Value = -0.08069646398190604
$ws2.Cells.Item(21, 3).Value = 0.07089796134946129
$ws2.Cells.Item(21, 4).Value = 0.1345297551567557
$ws2.Cells.Item(21, 5).Value = 0.8476251059191149
$ws2.Cells.Item(21, 6).Value = -0.1136918292190382

# --- Sheet 3: All Refactor Variances ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range('I2').Value = 4.983165035304014
$ws3.Range('I3').Value = 0.2491582517652007
$ws3.Range('I4').Value = 0.2491582517652007
$ws3.Range('I5').Value = 0.3608673071387364
$ws3.Range('J2').Value = 3.692480325572514
$ws3.Range('J3').Value = 0.1846240162786257
$ws3.Range('J4').Value = 0.4337822680438264
$ws3.Range('J5').Value = 0.267399418303798
$ws3.Range('K2').Value = 1.902149211798767
$ws3.Range('K3').Value = 0.09510746058993837
$ws3.Range('K4').Value = 0.5288897286337648
$ws3.Range('K5').Value = 0.1377484909640392
$ws3.Range('L2').Value = 1.763128716214867
$ws3.Range('L3').Value = 0.08815643581074337
$ws3.Range('L4').Value = 0.6170461644445082
$ws3.Range('L5').Value = 0.127681003428902
$ws3.Range('M2').Value = 1.467933697393225
$ws3.Range('M3').Value = 0.07339668486966126
$ws3.Range('M4').Value = 0.6904428493141694
$ws3.Range('M5').Value = 0.1063037801645243

# --- Sheet 4: Final Variances ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range('B2').Value = 4.983165035304014
$ws4.Range('B3').Value = 0.2491582517652007
$ws4.Range('B4').Value = 0.2491582517652007
$ws4.Range('B5').Value = 0.3608673071387364
$ws4.Range('C2').Value = 3.692480325572514
$ws4.Range('C3').Value = 0.1846240162786257
$ws4.Range('C4').Value = 0.4337822680438264
$ws4.Range('C5').Value = 0.267399418303798
$ws4.Range('D2').Value = 1.902149211798767
$ws4.Range('D3').Value = 0.09510746058993837
$ws4.Range('D4').Value = 0.5288897286337648
$ws4.Range('D5').Value = 0.1377484909640392
$ws4.Range('E2').Value = 1.763128716214867
$ws4.Range('E3').Value = 0.08815643581074337
$ws4.Range('E4').Value = 0.6170461644445082
$ws4.Range('E5').Value = 0.127681003428902
$ws4.Range('F2').Value = 1.467933697393225
$ws4.Range('F3').Value = 0.07339668486966126
$ws4.Range('F4').Value = 0.6904428493141694
$ws4.Range('F5').Value = 0.1063037801645243

# --- Sheet 5: Included and Excluded ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range('B2').Value = '[[''QEXTRCT'', ''QEDLESHI'', ''QNOHLTH'', ''QPOVTY'', ''QESL'', ''QHISPC'', ''PPUNIT'', ''PERCAP'', ''QFAM'', ''QFHH'', ''QBLACK'', ''QSERV'', ''QRICH'', ''QFEMLBR'', ''QFEMALE'', ''QAGEDEP'', ''QRENTER'', ''MEDAGE'', ''QSSBEN'', ''QNOAUTO'']]'
